# Applies the "Normalized noise level values for the test set" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Clean up the leftover "useless" number-format style (old cellXfs
#    index 1 -> numFmtId 0 / applyNumberFormat) on the cells that had it
#    in the first table (rows 14-19) and in the third table (rows 26-30).
#    Resetting the style to "Normal" removes the explicit style index,
#    matching cells that should be unstyled in the target workbook.
# ---------------------------------------------------------------------
$ws.Range("B14:K19").Style = "Normal"
$ws.Range("G26:K26").Style = "Normal"
$ws.Range("C27").Style = "Normal"
$ws.Range("C28").Style = "Normal"
$ws.Range("C29").Style = "Normal"
$ws.Range("C30").Style = "Normal"

# C16 keeps the custom "0.0000" number format (it just moves from the old
# unused cellXfs slot to the remaining one).
$ws.Range("C16").NumberFormat = "0.0000"

# ---------------------------------------------------------------------
# 2. Remove the empty, styled spacer row 20 entirely (not a row shift -
#    rows below keep their row numbers).
# ---------------------------------------------------------------------
$ws.Range("A20:K20").Clear()

# ---------------------------------------------------------------------
# 3. Add the new "normalized" table title, header and data (rows 33-40).
# ---------------------------------------------------------------------
$ws.Range("A33").Value = "Using mean squared error - rgb image - color channels separated - normalized (values between 0 and 5)"

$ws.Range("B35").Value = "image 1 "
$ws.Range("C35").Value = "image 2"
$ws.Range("D35").Value = "image 3"
$ws.Range("E35").Value = "image 4"
$ws.Range("F35").Value = "image 5"
$ws.Range("G35").Value = "image 6"
$ws.Range("H35").Value = "image 7"
$ws.Range("I35").Value = "image 8"
$ws.Range("J35").Value = "image 9"
$ws.Range("K35").Value = "image 10"

$ws.Range("A36").Value = "noise level 1"
$ws.Range("B36").Value = 0.16
$ws.Range("C36").Value = 0.02
$ws.Range("D36").Value = 0.13
$ws.Range("E36").Value = 0.09
$ws.Range("F36").Value = 0.02
$ws.Range("G36").Value = 0.06
$ws.Range("H36").Value = 0.15
$ws.Range("I36").Value = 0.11
$ws.Range("J36").Value = 0.05
$ws.Range("K36").Value = 0.04

$ws.Range("A37").Value = "noise level 2"
$ws.Range("B37").Value = 0.36
$ws.Range("C37").Value = 0.14
$ws.Range("D37").Value = 0.45
$ws.Range("E37").Value = 0.56
$ws.Range("F37").Value = 0.19
$ws.Range("G37").Value = 0.33
$ws.Range("H37").Value = 0.88
$ws.Range("I37").Value = 0.69
$ws.Range("J37").Value = 0.37
$ws.Range("K37").Value = 0.24

$ws.Range("A38").Value = "noise level 3"
$ws.Range("B38").Value = 0.8
$ws.Range("C38").Value = 0.52
$ws.Range("D38").Value = 1.19
$ws.Range("E38").Value = 1.46
$ws.Range("F38").Value = 0.64
$ws.Range("G38").Value = 0.75
$ws.Range("H38").Value = 2.09
$ws.Range("I38").Value = 1.75
$ws.Range("J38").Value = 0.97
$ws.Range("K38").Value = 0.79

$ws.Range("A39").Value = "noise level 4"
$ws.Range("B39").Value = 1.4
$ws.Range("C39").Value = 0.98
$ws.Range("D39").Value = 2.27
$ws.Range("E39").Value = 2.68
$ws.Range("F39").Value = 1.26
$ws.Range("G39").Value = 1.29
$ws.Range("H39").Value = 3.44
$ws.Range("I39").Value = 3.17
$ws.Range("J39").Value = 1.89
$ws.Range("K39").Value = 1.64

$ws.Range("A40").Value = "noise level 5"
$ws.Range("B40").Value = 2.24
$ws.Range("C40").Value = 1.56
$ws.Range("D40").Value = 3.74
$ws.Range("E40").Value = 4.31
$ws.Range("F40").Value = 2.1
$ws.Range("G40").Value = 2.03
$ws.Range("H40").Value = 4.94
$ws.Range("I40").Value = 5
$ws.Range("J40").Value = 3.16
$ws.Range("K40").Value = 2.85

# ---------------------------------------------------------------------
# 4. Update the view: scrolled down to row 10 and the active selection
#    moved to M34 (a cell just past the new table).
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M34").Select()
